$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new numeric value, per the authoritative diff
$updates = @{
    "E7" = 185
    "E8" = 40
    "E13" = 46
    "E14" = 8
    "E15" = 124
    "E16" = 30
    "E18" = 33
    "E19" = 7
    "E21" = 142
    "E22" = 38
    "E26" = 87
    "E27" = 15
    "E32" = 614
    "E33" = 98
    "E38" = 487
    "E39" = 143
    "E41" = 888
    "E42" = 231
    "E44" = 998
    "E45" = 294
    "E47" = 1102
    "E48" = 332
    "E50" = 1160
    "E51" = 377
    "E53" = 1136
    "E54" = 365
    "E56" = 1115
    "E57" = 443
    "E59" = 654
    "E60" = 243
    "E62" = 1250
    "E63" = 546
    "E68" = 775
    "E69" = 329
    "E71" = 925
    "E72" = 512
    "E74" = 918
    "E75" = 554
    "E77" = 707
    "E78" = 622
    "E80" = 850
    "E81" = 847
    "E83" = 672
    "E84" = 632
    "E86" = 967
    "E87" = 497
    "E89" = 700
    "E90" = 450
    "E92" = 603
    "E93" = 512
    "E95" = 704
    "E96" = 495
    "E98" = 850
    "E99" = 642
    "E101" = 914
    "E102" = 619
    "E104" = 339
    "E105" = 399
    "E107" = 744
    "E108" = 1028
    "E109" = 51
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
